# Apply the commit's data changes to the "Pasos" worksheet:
#  - D11 text is updated to clarify that "npm run" must be used before
#    generating the Home component.
#  - A new row (14) is inserted documenting that after "git pull" the
#    project must be recompiled with "npm install".
#  - The active selection is left on D13, matching the saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of instructions (npm install / recompile after git pull).
$ws.Range("D14").Value = "npm install"
$ws.Range("E14").Value = "recompilar despues de hacer pull en git"

# Update existing instruction cell (D11) with the corrected text.
$ws.Range("D11").Value = '"npm run" ng generate component Home'

# Restore the selection to the cell that was active when the file was saved.
$ws.Range("D13").Select()
